$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header K1: "Copia de NumeroTarjeta" -> "CopiadeNumeroTarjeta"
$ws.Range("K1").Value = "CopiadeNumeroTarjeta"

# Data row2 J2: "RONALD*RODRIGUEZ" -> "NOMBRE1*APELLIDO1"
$ws.Range("J2").Value = "NOMBRE1*APELLIDO1"

# New header column L1: "ID_BPM"
$ws.Range("L1").Value = "ID_BPM"

# Match the formatting of the new header cell to the existing header style (K1)
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

# Update selection to L1
$ws.Range("L1").Select()
